$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the seat number,
# shifting everything (STUDENT_NAME, MOTHER NAME, SUB1..SUB5, data rows) one column right.
$ws.Columns("A:A").Insert()

# Header for the new seat-number column
$ws.Range("A4").Value = "SEAT NO:"

# Seat numbers for each student row
$ws.Range("A5").Value = "s1"
$ws.Range("A6").Value = "s2"
$ws.Range("A7").Value = "s3"
$ws.Range("A8").Value = "s4"
$ws.Range("A9").Value = "s5"
$ws.Range("A10").Value = "s6"
$ws.Range("A11").Value = "s7"
$ws.Range("A12").Value = "s8"

# Match the final selection recorded in the saved file
$ws.Range("C9").Select()
